$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.897.90'
$ws.Range("E2").Value = '  -2.93%  '
$ws.Range("D3").Value = '2.654.71'
$ws.Range("E3").Value = '  -1.11%  '
$ws.Range("E4").Value = '  +0.15%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '523.22'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.33%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '144.36'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -1.32%  '
$ws.Range("E7").Value = '  +0.25%  '
$ws.Range("E8").Value = '  -1.44%  '
$ws.Range("E9").Value = '  +8.41%  '
$ws.Range("E10").Value = '  -2.82%  '
$ws.Range("E11").Value = '  -1.85%  '
$ws.Range("E12").Value = '  +1.32%  '
$ws.Range("D13").Value = '3.119.51'
$ws.Range("E13").Value = '  -0.95%  '
$ws.Range("D14").Value = '58.900.42'
$ws.Range("E14").Value = '  -2.58%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '21.05'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -1.05%  '
$ws.Range("B16").Value = 'ShibaInu'
$ws.Range("C16").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000136'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -1.88%  '
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '2.651.63'
$ws.Range("E17").Value = '  -3.22%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '338.87'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -3.42%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.37'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -4.13%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.37'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -1.77%  '
$ws.Range("E21").Value = '  +0.47%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.999'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +0.32%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '63.77'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +1.06%  '
$ws.Range("E24").Value = '  -1.31%  '
$ws.Range("E25").Value = '  -1.58%  '
$ws.Range("E26").Value = '  +0.89%  '
$ws.Range("E27").Value = '  -1.57%  '
$ws.Range("E28").Value = '  -2.66%  '
$ws.Range("E29").Value = '  -3.26%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.998'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +0.06%  '
$ws.Range("E31").Value = '  -0.18%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '18.84'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -1.26%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '149.69'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +0.66%  '
$ws.Range("E34").Value = '  -4.78%  '
$ws.Range("E35").Value = '  -3.14%  '
$ws.Range("E36").Value = '  -6.72%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.869'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -0.37%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '36.77'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -0.06%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.46'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -6.61%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.59'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -2.93%  '
$ws.Range("E41").Value = '  +0.73%  '
$ws.Range("E42").Value = '  +0.14%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '19.89'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -1.08%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '275.40'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -2.49%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0968'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -2.19%  '
$ws.Range("E46").Value = '  +2.07%  '
$ws.Range("B47").Value = 'Maker'
$ws.Range("C47").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D47").Value = '2.046.37'
$ws.Range("E47").Value = '  -4.06%  '
$ws.Range("B48").Value = 'Hedera'
$ws.Range("C48").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0530'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -1.98%  '
$ws.Range("E49").Value = '  -3.09%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '18.91'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -0.95%  '
$ws.Range("E51").Value = '  -3.10%  '
